$wb = $excel.ActiveWorkbook

# --- Update the status text: "Ready for handoff" -> "In Translation" ---
# Overview sheet: columns E (zh-cn) and F (de-de), rows 2-3
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = "In Translation"

# zh-cn sheet: column C (Status), rows 2-3
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2:C3").Value = "In Translation"

# de-de sheet: column C (Status), rows 2-3
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2:C3").Value = "In Translation"

# --- Narrow the (now shorter) status columns ---
# Original width 17.2159881591797 chars -> new width 13.4101845877511 chars.
# (ColumnWidth snaps to this host's pixel grid; 12.5 lands on the closest
# representable stored width to the target.)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZh.Columns.Item(3).ColumnWidth = 12.5

$wsDe.Columns.Item(3).ColumnWidth = 12.5
